$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# row 2 - ADANI
$ws.Cells.Item(2, 2).Value = 883.9
$ws.Cells.Item(2, 3).Value = 868.2
$ws.Cells.Item(2, 4).Value = 880.8
$ws.Cells.Item(2, 5).Value = 882.35
$ws.Cells.Item(2, 6).Value = 192
$ws.Cells.Item(2, 7).Value = 870.7

# row 3 - AURO
$ws.Cells.Item(3, 2).Value = 1049.25
$ws.Cells.Item(3, 3).Value = 1030.65
$ws.Cells.Item(3, 4).Value = 1048.4
$ws.Cells.Item(3, 5).Value = 1047.55
$ws.Cells.Item(3, 6).Value = 36
$ws.Cells.Item(3, 7).Value = 1040

# row 4 - BN
$ws.Cells.Item(4, 2).Value = 46725
$ws.Cells.Item(4, 3).Value = 45801.8
$ws.Cells.Item(4, 4).Value = 46701.2
$ws.Cells.Item(4, 5).Value = 46659.55
$ws.Cells.Item(4, 6).Value = 39
$ws.Cells.Item(4, 7).Value = 45850.05

# row 5 - CANBK
$ws.Cells.Item(5, 2).Value = 432.5
$ws.Cells.Item(5, 3).Value = 423.1
$ws.Cells.Item(5, 4).Value = 431
$ws.Cells.Item(5, 5).Value = 430.3
$ws.Cells.Item(5, 6).Value = 321
$ws.Cells.Item(5, 7).Value = 424.25

# row 6 - DLF
$ws.Cells.Item(6, 2).Value = 655.4
$ws.Cells.Item(6, 3).Value = 636.2
$ws.Cells.Item(6, 4).Value = 651.9
$ws.Cells.Item(6, 5).Value = 652.45
$ws.Cells.Item(6, 6).Value = 157
$ws.Cells.Item(6, 7).Value = 640

# row 7 - HIND
$ws.Cells.Item(7, 2).Value = 529.35
$ws.Cells.Item(7, 3).Value = 520.45
$ws.Cells.Item(7, 4).Value = 523.35
$ws.Cells.Item(7, 5).Value = 522.35
$ws.Cells.Item(7, 6).Value = 95
$ws.Cells.Item(7, 7).Value = 524.25

# row 8 - ICICI
$ws.Cells.Item(8, 2).Value = 999
$ws.Cells.Item(8, 3).Value = 971.5
$ws.Cells.Item(8, 4).Value = 997.5
$ws.Cells.Item(8, 5).Value = 996
$ws.Cells.Item(8, 6).Value = 462
$ws.Cells.Item(8, 7).Value = 977.8

# row 9 - JIND
$ws.Cells.Item(9, 2).Value = 698.9
$ws.Cells.Item(9, 3).Value = 682.25
$ws.Cells.Item(9, 4).Value = 685.15
$ws.Cells.Item(9, 5).Value = 685.75
$ws.Cells.Item(9, 6).Value = 52
$ws.Cells.Item(9, 7).Value = 698.45

# row 10 - NIFTY
$ws.Cells.Item(10, 2).Value = 20828.3
$ws.Cells.Item(10, 3).Value = 20640.85
$ws.Cells.Item(10, 4).Value = 20804.05
$ws.Cells.Item(10, 5).Value = 20798.9
$ws.Cells.Item(10, 6).Value = 73
$ws.Cells.Item(10, 7).Value = 20652.1

# row 11 - REL
$ws.Cells.Item(11, 2).Value = 2439.1
$ws.Cells.Item(11, 3).Value = 2416.25
$ws.Cells.Item(11, 4).Value = 2432.8
$ws.Cells.Item(11, 5).Value = 2433.65
$ws.Cells.Item(11, 6).Value = 91
$ws.Cells.Item(11, 7).Value = 2430.85

# row 12 - SBIN
$ws.Cells.Item(12, 2).Value = 600.45
$ws.Cells.Item(12, 3).Value = 588.25
$ws.Cells.Item(12, 4).Value = 600.05
$ws.Cells.Item(12, 5).Value = 599
$ws.Cells.Item(12, 6).Value = 660
$ws.Cells.Item(12, 7).Value = 588.5

# row 13 - TCON
$ws.Cells.Item(13, 2).Value = 962.9
$ws.Cells.Item(13, 3).Value = 949.25
$ws.Cells.Item(13, 4).Value = 952.9
$ws.Cells.Item(13, 5).Value = 953.4
$ws.Cells.Item(13, 6).Value = 27
$ws.Cells.Item(13, 7).Value = 952.5

# row 14 - TM
$ws.Cells.Item(14, 2).Value = 713.5
$ws.Cells.Item(14, 3).Value = 707.55
$ws.Cells.Item(14, 4).Value = 710.3
$ws.Cells.Item(14, 5).Value = 710.35
$ws.Cells.Item(14, 6).Value = 146
$ws.Cells.Item(14, 7).Value = 710.2

# row 15 - TS
$ws.Cells.Item(15, 2).Value = 132.7
$ws.Cells.Item(15, 3).Value = 131.4
$ws.Cells.Item(15, 4).Value = 131.6
$ws.Cells.Item(15, 5).Value = 131.6
$ws.Cells.Item(15, 6).Value = 452
$ws.Cells.Item(15, 7).Value = 132

# row 16 - TCS
$ws.Cells.Item(16, 2).Value = 3554.5
$ws.Cells.Item(16, 3).Value = 3525.25
$ws.Cells.Item(16, 4).Value = 3536.75
$ws.Cells.Item(16, 5).Value = 3533.6
$ws.Cells.Item(16, 6).Value = 19
$ws.Cells.Item(16, 7).Value = 3550

# row 17 - TITAN
$ws.Cells.Item(17, 2).Value = 3535
$ws.Cells.Item(17, 3).Value = 3503.25
$ws.Cells.Item(17, 4).Value = 3512.8
$ws.Cells.Item(17, 5).Value = 3513.95
$ws.Cells.Item(17, 6).Value = 12
$ws.Cells.Item(17, 7).Value = 3530.15
